$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a header row (row 1) plus data rows, columns A:E, where
# column D holds "codeforiati:group-name" values and column E holds
# "codeforiati:group-code" values (including the D1/E1 header labels
# themselves). The source data had the group-name / group-code values
# transposed; this swaps D and E for every row (header + data) so each
# column holds the correct values again.

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $dCell.Value2 = $eVal
    $eCell.Value2 = $dVal
}
